# Regenerate the K (strikeouts) column values (column G) for the 2023
# wilson_bryse save_data sheet. The original sheet used a different
# "Strike#" derived count; this recalculates/rewrites the true K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 4
    15 = 3
    16 = 0
    17 = 2
    18 = 0
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 5
    25 = 0
    26 = 0
    27 = 2
    28 = 2
    29 = 1
    30 = 1
    31 = 1
    32 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 1
    44 = 2
    45 = 2
    46 = 1
    47 = 0
    48 = 1
    49 = 4
    50 = 3
    51 = 1
    52 = 2
    53 = 3
    54 = 2
    55 = 1
    56 = 2
    57 = 1
    58 = 4
    59 = 2
    60 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
